# Add "weight" (E) and "main_score" (F) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Copy the formatting of the existing header cell (bold / centered /
# bordered, style index 1) onto the two new header cells so we reuse the
# existing style instead of creating a brand-new one.
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E1").Value = "weight"
$ws.Range("F1").Value = "main_score"

# --- Data rows ----------------------------------------------------------
# Each contiguous block of rows shares the same sentence group, and thus
# the same weight / main_score pair.
$blocks = @(
    @{ Start = 2;  End = 8;  Weight = 40; Score = 0.5336644649505615 },
    @{ Start = 9;  End = 18; Weight = 40; Score = 0.4905979931354523 },
    @{ Start = 19; End = 38; Weight = 40; Score = 0.5350711345672607 },
    @{ Start = 39; End = 40; Weight = 40; Score = 0.4114404618740082 },
    @{ Start = 41; End = 47; Weight = 40; Score = 0.5032179951667786 },
    @{ Start = 48; End = 48; Weight = 40; Score = 0.4224454164505005 },
    @{ Start = 49; End = 59; Weight = 20; Score = 0.6599290370941162 },
    @{ Start = 60; End = 70; Weight = 5;  Score = 0.4740377962589264 },
    @{ Start = 71; End = 73; Weight = 30; Score = 0.4175359904766083 },
    @{ Start = 74; End = 74; Weight = 30; Score = 0.4401599764823914 },
    @{ Start = 75; End = 88; Weight = 5;  Score = 0.5501365065574646 }
)

foreach ($block in $blocks) {
    $ws.Range("E$($block.Start):E$($block.End)").Value = $block.Weight
    $ws.Range("F$($block.Start):F$($block.End)").Value = $block.Score
}
